$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy formatting from row 18 into the new row 19 before filling values in,
# so the new row matches the existing styling (s="6" for A/B/D/E, s="2" for C).
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)  # xlPasteFormats

# Existing TestCase_F17 row (18) used to be the one marked PASS; now that
# TestCase_F18 is the passing case, F17 becomes SKIP.
$ws.Range("E18").Value = "SKIP"

# New row 19: TestCase_F18 / OPQA-1099
$ws.Range("A19").Value = "TestCase_F18"
$ws.Range("B19").Value = "OPQA-1099"
$ws.Range("C19").Value = "Verify that Featured Post move down when new notification event occur"
$ws.Range("D19").Value = "Y"
$ws.Range("E19").Value = "PASS"

# Update selection / active cell to the newly added row, matching the saved view state.
$ws.Range("A19").Select()
